$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.08243498013860336
$ws.Range("J2").Value = 0.08243498013860337
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 0.1877172950347778
$ws.Range("R2").Value = 1.689455655313
$ws.Range("S2").Value = 0.0009093100199700476
$ws.Range("T2").Value = 0.0009093100199700476

$ws.Range("I3").Value = 0.08243498013860336
$ws.Range("J3").Value = 0.08243498013860337
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("S3").Value = 0.06163019849033417
$ws.Range("T3").Value = 0.06163019849033417

$ws.Range("I4").Value = 0.08243498013860336
$ws.Range("J4").Value = 0.08243498013860337
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 4.069421357249334
$ws.Range("R4").Value = 36.624792215244
$ws.Range("S4").Value = 0.01971243840340537
$ws.Range("T4").Value = 0.01971243840340537

$ws.Range("I5").Value = 0.08243498013860336
$ws.Range("J5").Value = 0.08243498013860337
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 0.03778524499233334
$ws.Range("R5").Value = 0.340067204931
$ws.Range("S5").Value = 0.0001830332248937759
$ws.Range("T5").Value = 0.0001830332248937759

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.524162666666667
$ws.Range("H6").Value = 4.572488
$ws.Range("I6").Value = 0.9175650198613966
$ws.Range("J6").Value = 0.9175650198613967
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 2.089438527883555
$ws.Range("R6").Value = 18.804946750952
$ws.Range("S6").Value = 0.01012132307342265
$ws.Range("T6").Value = 0.01012132307342265

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.524162666666667
$ws.Range("H7").Value = 4.572488
$ws.Range("I7").Value = 0.9175650198613966
$ws.Range("J7").Value = 0.9175650198613967
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("Q7").Value = 141.6156298498249
$ws.Range("R7").Value = 1274.540668648424
$ws.Range("S7").Value = 0.6859917259246563
$ws.Range("T7").Value = 0.6859917259246563

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.524162666666667
$ws.Range("H8").Value = 4.572488
$ws.Range("I8").Value = 0.9175650198613966
$ws.Range("J8").Value = 0.9175650198613967
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 45.29580382273066
$ws.Range("R8").Value = 407.662234404576
$ws.Range("S8").Value = 0.2194146696551099
$ws.Range("T8").Value = 0.2194146696551099

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.524162666666667
$ws.Range("H9").Value = 4.572488
$ws.Range("I9").Value = 0.9175650198613966
$ws.Range("J9").Value = 0.9175650198613967
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 0.4205789704026666
$ws.Range("R9").Value = 3.785210733624
$ws.Range("S9").Value = 0.002037301208207683
$ws.Range("T9").Value = 0.002037301208207682

